$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.6545652718822623
$ws.Range("C2").Value = 2919.202174992006
$ws.Range("D2").Value = 18.71679738969934
$ws.Range("E2").Value = 13.86384647080068
$ws.Range("G2").Value = 2952.437384124388
